$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164, shifting rows 164:181 down to 165:182.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with its data.
$ws.Cells.Item(164, 1).Value = 11
$ws.Cells.Item(164, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(164, 3).Value = "Bíobío"
$ws.Cells.Item(164, 4).Value = 45154
$ws.Cells.Item(164, 5).Value = 8
$ws.Cells.Item(164, 6).Value = 100112001
$ws.Cells.Item(164, 7).Value = "Berenjena"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 230
$ws.Cells.Item(164, 11).Value = 7000
$ws.Cells.Item(164, 12).Value = 8000
$ws.Cells.Item(164, 13).Value = 7652
$ws.Cells.Item(164, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(164, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(164, 16).Value = 128
$ws.Cells.Item(164, 17).Value = 60
$ws.Cells.Item(164, 18).Value = "Hortaliza"
